# "Generate Report for Archive"
#
# The handback status "Ready for handoff" is now "In Translation" for the
# zh-cn / de-de locales (and their roll-up on the Overview sheet), and the
# narrower text means the Status columns can be narrower too.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: the zh-cn / de-de status roll-up columns (E & F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Narrow the E:F columns to fit the new, shorter status text.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status column (C) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.5
